$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: becomes old Row 4's data (Especial, date 45008, etc.)
$ws.Range("D3").Value = 45008
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("S3").Value = 3500

# Row 4: becomes old Row 5's data (Primera, date 45008 unchanged)
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000

# Row 5: becomes old Row 3's data (date 44995, Primera unchanged)
$ws.Range("D5").Value = 44995
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5500
$ws.Range("P5").Value = 5750
$ws.Range("S5").Value = 2875
